$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.529.80"
$ws.Range("E2").Value = "  -1.53%  "
$ws.Range("D3").Value = "1.670.64"
$ws.Range("E3").Value = "  -1.92%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "'312.74"
$ws.Range("E5").Value = "  -0.95%  "
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").Value = "'0.3907"
$ws.Range("E7").Value = "  -4.01%  "
$ws.Range("D8").Value = "'0.3940"
$ws.Range("E8").Value = "  -3.03%  "
$ws.Range("D9").Value = "'1.001"
$ws.Range("E9").Value = "  -0.20%  "
$ws.Range("D10").Value = "'51.35"
$ws.Range("E10").Value = "  -4.50%  "
$ws.Range("D11").Value = "'1.401"
$ws.Range("E11").Value = "  -4.46%  "
$ws.Range("D12").Value = "'0.08633"
$ws.Range("E12").Value = "  -2.07%  "
$ws.Range("D13").Value = "'25.30"
$ws.Range("E13").Value = "  -1.82%  "
$ws.Range("D14").Value = "'7.330"
$ws.Range("E14").Value = "  -2.11%  "
$ws.Range("D15").Value = "'0.00001323"
$ws.Range("E15").Value = "  -2.10%  "
$ws.Range("D16").Value = "'7.733"
$ws.Range("E16").Value = "  -3.87%  "
$ws.Range("D17").Value = "1.672.41"
$ws.Range("E17").Value = "  -3.00%  "
$ws.Range("D18").Value = "'93.30"
$ws.Range("D19").Value = "'0.07020"
$ws.Range("E19").Value = "  -2.09%  "
$ws.Range("D20").Value = "'21.05"
$ws.Range("E20").Value = "  +0.48%  "
$ws.Range("D21").Value = "'7.052"
$ws.Range("E21").Value = "  -2.50%  "
$ws.Range("D22").Value = "'1.002"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("E23").Value = "  -4.52%  "
$ws.Range("D24").Value = "24.534.02"
$ws.Range("E24").Value = "  -1.51%  "
$ws.Range("D25").Value = "'2.357"
$ws.Range("E25").Value = "  +1.32%  "
$ws.Range("D26").Value = "'2.744"
$ws.Range("E26").Value = "  -4.75%  "
$ws.Range("D27").Value = "'23.15"
$ws.Range("E27").Value = "  +0.45%  "
$ws.Range("D28").Value = "'5.859"
$ws.Range("E28").Value = "  -13.35%  "
$ws.Range("D29").Value = "'160.27"
$ws.Range("E29").Value = "  -2.70%  "
$ws.Range("D30").Value = "'146.28"
$ws.Range("E30").Value = "  +0.96%  "
$ws.Range("D31").Value = "'8.375"
$ws.Range("E31").Value = "  +1.67%  "
$ws.Range("D32").Value = "'2.507"
$ws.Range("E32").Value = "  +10.51%  "
$ws.Range("D33").Value = "1.857.61"
$ws.Range("E33").Value = "  -3.78%  "
$ws.Range("E34").Value = "  -4.96%  "
$ws.Range("D35").Value = "'7.004"
$ws.Range("E35").Value = "  -4.06%  "
$ws.Range("D36").Value = "'0.03032"
$ws.Range("E36").Value = "  -5.22%  "
$ws.Range("D37").Value = "'0.2809"
$ws.Range("E37").Value = "  -1.19%  "
$ws.Range("D38").Value = "'0.9896"
$ws.Range("E38").Value = "  -2.78%  "
$ws.Range("D39").Value = "'0.09467"
$ws.Range("E39").Value = "  +0.30%  "
$ws.Range("D40").Value = "'1.512"
$ws.Range("E40").Value = "  +2.95%  "
$ws.Range("E41").Value = "  -6.02%  "
$ws.Range("D42").Value = "'0.7905"
$ws.Range("E42").Value = "  -7.01%  "
$ws.Range("E43").Value = "  -3.35%  "
$ws.Range("D44").Value = "'16.42"
$ws.Range("E44").Value = "  -7.34%  "
$ws.Range("D45").Value = "'0.7112"
$ws.Range("E45").Value = "  -4.27%  "
$ws.Range("D46").Value = "'2.546"
$ws.Range("E46").Value = "  -6.58%  "
$ws.Range("D47").Value = "'4.182"
$ws.Range("E47").Value = "  -1.23%  "
$ws.Range("D48").Value = "'0.08619"
$ws.Range("E48").Value = "  +3.15%  "
$ws.Range("E49").Value = "  -0.16%  "
$ws.Range("D50").Value = "'1.324"
$ws.Range("E50").Value = "  -5.40%  "
$ws.Range("D51").Value = "'137.39"
$ws.Range("E51").Value = "  -3.29%  "
